{"js": "// Append three new bulleted paragraphs right after the\n// \"Resolution des erreurs au demarrage de tomcat...\" list item, at the\n// end of the document body (before the section break), mirroring the\n// \"semaine 3 v3.6\" update: a new level-0 item describing the versioning\n// tool improvement, followed by two level-1 sub-items describing the\n// two proposed \"logiques\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText = \"Resolution des erreurs au demarrage de tomcat. Ex : log4j qui n'a pas acces au fichier de log\";\n\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === anchorText) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\nif (!anchor) {\n  throw new Error(\"Could not locate the anchor paragraph ('Resolution des erreurs au demarrage de tomcat...').\");\n}\n\n// New top-level (ilvl 0) bullet: reuses the anchor paragraph's list\n// formatting (same style/numId/ilvl) because it's inserted as a sibling\n// right after it.\nconst p1 = anchor.insertParagraph(\n  \"Ameliorer l'outil de versionning de documents.\",\n  \"After\"\n);\n// Second run appended to the same paragraph.\np1.getRange(\"End\").insertText(\n  \" 2 logiques de versionning propos\u00e9, a moi de choisir l'une d'elle. \",\n  \"End\"\n);\n\n// Sub-bullet (ilvl 1) #1.\nconst p2 = p1.insertParagraph(\n  \"Logique 1 : Comparer le nouveau document (avant insertion), avec le document actif. Puis ajouter dans la base de versionning, que les champs qui diff\u00e8re entre les 2. Et inserer le nouveau documents comme actif.\",\n  \"After\"\n);\np2.listItem.level = 1;\n\n// Sub-bullet (ilvl 1) #2.\nconst p3 = p2.insertParagraph(\n  \"Logique 2 : Inserer dans la base de versionning le document actif au complet, puis mettre le nouveau document comme actif. Et faire tourn\u00e9e une routine qui va comparer les champs et supprimer ceux qui sont identiques.\",\n  \"After\"\n);\np3.listItem.level = 1;\n\nawait context.sync();\n", "ps1": "# Append three new bulleted paragraphs right after the\n# \"Resolution des erreurs au demarrage de tomcat...\" list item, at the\n# end of the document body (before the section break), mirroring the\n# \"semaine 3 v3.6\" update: a new level-0 item describing the versioning\n# tool improvement, followed by two level-1 sub-items describing the\n# two proposed \"logiques\".\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n\n$target = $null\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $p = $paras.Item($i)\n    if ($p.Range.Text -like \"*Resolution des erreurs au demarrage de tomcat*\") {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Could not locate the anchor paragraph ('Resolution des erreurs au demarrage de tomcat...').\"\n}\n\n# New top-level (ilvl 0) bullet: InsertParagraphAfter clones the anchor\n# paragraph's list formatting (same style/numId/ilvl) since it becomes its\n# sibling right after it.\n$target.Range.InsertParagraphAfter()\n$p1 = $d.Paragraphs.Item($target.Index + 1)\n$p1.Range.InsertAfter(\"Ameliorer l'outil de versionning de documents.\")\n$p1.Range.InsertAfter(\" 2 logiques de versionning propos\u00e9, a moi de choisir l'une d'elle. \")\n\n# Sub-bullet (ilvl 1 / COM ListLevelNumber 2) #1.\n$p1.Range.InsertParagraphAfter()\n$p2 = $d.Paragraphs.Item($p1.Index + 1)\n$p2.Range.InsertAfter(\"Logique 1 : Comparer le nouveau document (avant insertion), avec le document actif. Puis ajouter dans la base de versionning, que les champs qui diff\u00e8re entre les 2. Et inserer le nouveau documents comme actif.\")\n$p2.Range.ListFormat.ListLevelNumber = 2\n\n# Sub-bullet (ilvl 1 / COM ListLevelNumber 2) #2.\n$p2.Range.InsertParagraphAfter()\n$p3 = $d.Paragraphs.Item($p2.Index + 1)\n$p3.Range.InsertAfter(\"Logique 2 : Inserer dans la base de versionning le document actif au complet, puis mettre le nouveau document comme actif. Et faire tourn\u00e9e une routine qui va comparer les champs et supprimer ceux qui sont identiques.\")\n$p3.Range.ListFormat.ListLevelNumber = 2\n"}
